$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 9 new data rows (rows 29-37), matching the target diff.
# Columns C (Tel) and G (Cep) contain digit-only strings that must stay
# as text (preserving values such as a leading zero in the CEP column),
# so their number format is forced to text ("@") before the value is set.

# Row 29
$ws.Range("A29").Value = 28
$ws.Range("B29").Value = 'joão P Santos'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = '11912345678'
$ws.Range("D29").Value = 25
$ws.Range("E29").Value = 'dsdh@gmm.com.br'
$ws.Range("F29").Value = 'Masculino'
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '06666000'
$ws.Range("H29").Value = 'Mario Veloso Serqueira'
$ws.Range("I29").Value = 43
$ws.Range("J29").Value = 'b'
$ws.Range("K29").Value = 'Carlos Drummord Andrade'
$ws.Range("L29").Value = 'Caracas'
$ws.Range("M29").Value = 'RN'

# Row 30
$ws.Range("A30").Value = 29
$ws.Range("B30").Value = 'joão P Santos'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = '11912345678'
$ws.Range("D30").Value = 25
$ws.Range("E30").Value = 'dsdh@gmm.com.br'
$ws.Range("F30").Value = 'Masculino'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '06666000'
$ws.Range("H30").Value = 'Mario Veloso Serqueira'
$ws.Range("I30").Value = 43
$ws.Range("J30").Value = 'b'
$ws.Range("K30").Value = 'Carlos Drummord Andrade'
$ws.Range("L30").Value = 'Caracas'
$ws.Range("M30").Value = 'RN'

# Row 31
$ws.Range("A31").Value = 30
$ws.Range("B31").Value = 'asdasdasdasdasdasd'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = '11912345678'
$ws.Range("D31").Value = 25
$ws.Range("E31").Value = 'dsdh@gmm.com.br'
$ws.Range("F31").Value = 'Masculino'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '06666000'
$ws.Range("H31").Value = 'Mario Veloso Serqueira'
$ws.Range("I31").Value = 43
$ws.Range("J31").Value = 'b'
$ws.Range("K31").Value = 'Carlos Drummord Andrade'
$ws.Range("L31").Value = 'Caracas'
$ws.Range("M31").Value = 'RN'

# Row 32
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = 'joão P Santos'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = '11912345678'
$ws.Range("D32").Value = 25
$ws.Range("E32").Value = 'dsdh@gmm.com.br'
$ws.Range("F32").Value = 'Feminino'
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '06666000'
$ws.Range("H32").Value = 'Mario Veloso Serqueira'
$ws.Range("I32").Value = 43
$ws.Range("J32").Value = 'b'
$ws.Range("K32").Value = 'Carlos Drummord Andrade'
$ws.Range("L32").Value = 'Caracas'
$ws.Range("M32").Value = 'RN'

# Row 33
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = 'Joao de souza cruz'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = '1961226556'
$ws.Range("D33").Value = 25
$ws.Range("E33").Value = 'dsdh@gmm.com.br'
$ws.Range("F33").Value = 'Feminino'
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '06666000'
$ws.Range("H33").Value = 'Mario Veloso Serqueira'
$ws.Range("I33").Value = 43
$ws.Range("J33").Value = 'b'
$ws.Range("K33").Value = 'Carlos Drummord Andrade'
$ws.Range("L33").Value = 'Caracas'
$ws.Range("M33").Value = 'RN'

# Row 34
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = 'joão P Santos'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = '1191234567'
$ws.Range("D34").Value = 25
$ws.Range("E34").Value = 'dsdh@gmm.com.br'
$ws.Range("F34").Value = 'Masculino'
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '06666000'
$ws.Range("H34").Value = 'Mario Veloso Serqueira'
$ws.Range("I34").Value = 43
$ws.Range("J34").Value = 'b'
$ws.Range("K34").Value = 'Carlos Drummord Andrade'
$ws.Range("L34").Value = 'Caracas'
$ws.Range("M34").Value = 'RN'

# Row 35
$ws.Range("A35").Value = 34
$ws.Range("B35").Value = 'joão P Santos'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = '1191234567'
$ws.Range("D35").Value = 25
$ws.Range("E35").Value = 'dsdh@gmm.com.br'
$ws.Range("F35").Value = 'Masculino'
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '06666000'
$ws.Range("H35").Value = 'Mario Veloso Serqueira'
$ws.Range("I35").Value = 43
$ws.Range("J35").Value = 'b'
$ws.Range("K35").Value = 'Carlos Drummord Andrade'
$ws.Range("L35").Value = 'Caracas'
$ws.Range("M35").Value = 'RN'

# Row 36
$ws.Range("A36").Value = 35
$ws.Range("B36").Value = 'joão P Santos'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = '11912345678'
$ws.Range("D36").Value = 25
$ws.Range("E36").Value = 'dsdh@gmm.com.br'
$ws.Range("F36").Value = 'Masculino'
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '06666000'
$ws.Range("H36").Value = 'Mario Veloso Serqueira'
$ws.Range("I36").Value = 43
$ws.Range("J36").Value = 'b'
$ws.Range("K36").Value = 'Carlos Drummord Andrade'
$ws.Range("L36").Value = 'Caracas'
$ws.Range("M36").Value = 'RN'

# Row 37
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = 'joão P Santos'
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = '11123456789'
$ws.Range("D37").Value = 25
$ws.Range("E37").Value = 'dsdh@gmm.com.br'
$ws.Range("F37").Value = 'Masculino'
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '06666000'
$ws.Range("H37").Value = 'Mario Veloso Serqueira'
$ws.Range("I37").Value = 43
$ws.Range("J37").Value = 'b'
$ws.Range("K37").Value = 'Carlos Drummord Andrade'
$ws.Range("L37").Value = 'Caracas'
$ws.Range("M37").Value = 'RN'

